# Update the crypto price/volume snapshot (Price + Volume(1h) columns),
# plus a couple of coin rows (41/42) that were re-ranked (name, link, price,
# volume all changed). NumberFormat is forced to Text ("@") before writing
# D/E values so numeric-looking strings (e.g. "243.93", "-0.06%") are stored
# verbatim as text instead of being parsed into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "243.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.03%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "14.02%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.148"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.07%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05679"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.41%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.531"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.89%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8467"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.44%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8562"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.02%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1333"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.23%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06907"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.25%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.02886"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.02%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09380"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.10%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001516"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.32%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04167"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-9.92%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006024"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-93.99%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006141"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.92%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-4.05%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.30%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.245"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.86%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.19%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03282"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.90%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.26%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.609"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.52%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.31%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001212"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.91%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004442"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-1.17%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "0.24%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03726"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.33%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005325"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-13.68%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1058"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.71%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.68%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009861"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "9.68%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005100"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.80%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.07%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.09990"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-30.61%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002799"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "21.06%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.07%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.07%"
